$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 34003.266
$ws.Range("J17").Value = 34003.266
$ws.Range("L17").Value = 102009.798
$ws.Range("N17").Value = -102345.798

$ws.Range("H112").Value = 2334.092
$ws.Range("I112").Value = 300
$ws.Range("J112").Value = 2432.1204
$ws.Range("K112").Value = 900
$ws.Range("L112").Value = 7296.361199999999
$ws.Range("M112").Value = 208
$ws.Range("N112").Value = -9512.361199999999

$ws.Range("H132").Value = 2023.3914
$ws.Range("I132").Value = 1517.9615
$ws.Range("K132").Value = 4553.8845
$ws.Range("M132").Value = -2023.8845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18339.162
$ws.Range("I32").Value = 20050.934
$ws.Range("J32").Value = 11003
$ws.Range("K32").Value = 20050.934
$ws.Range("L32").Value = 11003
$ws.Range("M32").Value = -19763.934
$ws.Range("N32").Value = -11577

$ws.Range("H41").Value = 3336.3
$ws.Range("I41").Value = 1984.7778
$ws.Range("J41").Value = 15500
$ws.Range("K41").Value = 1984.7778
$ws.Range("L41").Value = 15500
$ws.Range("M41").Value = -1570.7778
$ws.Range("N41").Value = -16328

$ws.Range("H61").Value = 3520
$ws.Range("I61").Value = 1989.3334
$ws.Range("K61").Value = 1989.3334
$ws.Range("M61").Value = -1777.3334

$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492

$ws.Range("H104").Value = 30000
$ws.Range("J104").Value = 30000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988

$ws.Range("H130").Value = 182895.67
$ws.Range("J130").Value = 182895.67
$ws.Range("L130").Value = 182895.67
$ws.Range("N130").Value = -192935.67

$ws.Range("H136").Value = 3520
$ws.Range("I136").Value = 1989.3334
$ws.Range("K136").Value = 5968.0002
$ws.Range("M136").Value = -3418.0002

$ws.Range("H139").Value = 54742.555
$ws.Range("J139").Value = 54256.824
$ws.Range("L139").Value = 54256.824
$ws.Range("N139").Value = -64536.824

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 2227.9092
$ws.Range("I37").Value = 331.4
$ws.Range("J37").Value = 3808.3333
$ws.Range("K37").Value = 331.4
$ws.Range("L37").Value = 3808.3333
$ws.Range("M37").Value = -194.4
$ws.Range("N37").Value = -4082.3333

$ws.Range("H134").Value = 2015.431
$ws.Range("I134").Value = 1674.6279
$ws.Range("K134").Value = 5023.8837
$ws.Range("M134").Value = -2488.8837

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3702.75
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3702.75
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3702.75
$ws.Range("N31").Value = -4292.75
$ws.Range("M31").ClearContents()

$ws.Range("H34").Value = 3702.75
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3702.75
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3702.75
$ws.Range("N34").Value = -4106.75
$ws.Range("M34").ClearContents()

$ws.Range("H58").Value = 1324918.8
$ws.Range("I58").Value = 2316872
$ws.Range("J58").Value = 2314.4167
$ws.Range("K58").Value = 2316872
$ws.Range("L58").Value = 2314.4167
$ws.Range("M58").Value = -2316669
$ws.Range("N58").Value = -2720.4167

$ws.Range("H87").Value = 44000
$ws.Range("J87").Value = 44000
$ws.Range("L87").Value = 44000
$ws.Range("N87").Value = -46372

$ws.Range("H90").Value = 44000
$ws.Range("J90").Value = 44000
$ws.Range("L90").Value = 132000
$ws.Range("N90").Value = -143856

$ws.Range("H134").Value = 1294.1522
$ws.Range("I134").Value = 1051.4324
$ws.Range("K134").Value = 3154.2972
$ws.Range("M134").Value = -619.2972

$ws.Range("H136").Value = 1324918.8
$ws.Range("I136").Value = 2316872
$ws.Range("J136").Value = 2314.4167
$ws.Range("K136").Value = 6950616
$ws.Range("L136").Value = 6943.250100000001
$ws.Range("M136").Value = -6948066
$ws.Range("N136").Value = -12043.2501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 33340106
$ws.Range("I3").Value = 90912260
$ws.Range("J3").Value = 8863
$ws.Range("K3").Value = 272736780
$ws.Range("L3").Value = 26589
$ws.Range("M3").Value = -272736668
$ws.Range("N3").Value = -26813

$ws.Range("H12").Value = 623381.9399999999
$ws.Range("I12").Value = 51.333332
$ws.Range("K12").Value = 153.999996
$ws.Range("M12").Value = 19.00000399999999

$ws.Range("H35").Value = 13250
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 13250
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 39750
$ws.Range("N35").Value = -40326
$ws.Range("M35").ClearContents()

$ws.Range("H88").Value = 50898.4
$ws.Range("J88").Value = 61869.5
$ws.Range("L88").Value = 185608.5
$ws.Range("N88").Value = -186464.5

$ws.Range("H91").Value = 50898.4
$ws.Range("J91").Value = 61869.5
$ws.Range("L91").Value = 185608.5
$ws.Range("N91").Value = -188572.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 216405.72
$ws.Range("I97").Value = 85806.664
$ws.Range("J97").Value = 1000000
$ws.Range("K97").Value = 85806.664
$ws.Range("L97").Value = 1000000
$ws.Range("M97").Value = -85310.664
$ws.Range("N97").Value = -1000992

$ws.Range("H123").Value = 10043.75
$ws.Range("J123").Value = 10043.75
$ws.Range("L123").Value = 10043.75
$ws.Range("N123").Value = -14943.75

$ws.Range("H131").Value = 36643.75
$ws.Range("J131").Value = 36643.75
$ws.Range("L131").Value = 36643.75
$ws.Range("N131").Value = -46723.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3262.5334
$ws.Range("I7").Value = 3094.8333
$ws.Range("K7").Value = 3094.8333
$ws.Range("M7").Value = -2982.8333

$ws.Range("H126").Value = 3262.5334
$ws.Range("I126").Value = 3094.8333
$ws.Range("K126").Value = 9284.499899999999
$ws.Range("M126").Value = -6814.499899999999

$ws.Range("H131").Value = 16384
$ws.Range("J131").Value = 16384
$ws.Range("L131").Value = 16384
$ws.Range("N131").Value = -26464

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 102568.84
$ws.Range("I81").Value = 93581.37
$ws.Range("J81").Value = 152000
$ws.Range("K81").Value = 187162.74
$ws.Range("L81").Value = 304000
$ws.Range("M81").Value = -186101.74
$ws.Range("N81").Value = -306122

$ws.Range("H84").Value = 102568.84
$ws.Range("I84").Value = 93581.37
$ws.Range("J84").Value = 152000
$ws.Range("K84").Value = 935813.7
$ws.Range("L84").Value = 1520000
$ws.Range("M84").Value = -930509.7
$ws.Range("N84").Value = -1530608

$ws.Range("H100").Value = 11404.842
$ws.Range("I100").Value = 17210.5
$ws.Range("J100").Value = 1452.2858
$ws.Range("K100").Value = 34421
$ws.Range("L100").Value = 2904.5716
$ws.Range("M100").Value = -33880
$ws.Range("N100").Value = -3986.5716

$ws.Range("H123").Value = 22272.924
$ws.Range("J123").Value = 22272.924
$ws.Range("L123").Value = 22272.924
$ws.Range("N123").Value = -32072.924

$ws.Range("H135").Value = 29800
$ws.Range("J135").Value = 29800
$ws.Range("L135").Value = 29800
$ws.Range("N135").Value = -39940

$ws.Range("H136").Value = 1359.1852
$ws.Range("I136").Value = 1373
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 4119
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -1569
$ws.Range("N136").Value = -8100
